$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Revert "User data 3.0": remove the "budget-type" column (column B) from the
# wide-format data sheet, shifting the year columns (2013-2016) one column
# to the left so the layout goes back to id, 2013, 2014, 2015, 2016.
$ws.Columns.Item(2).Delete()
